$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/1/2024  Through  4/7/2024"

# --- Donor cells for styles (unchanged styles used as format source) ---
$styleDonor14 = "A15"   # string style (t=s), General numFmt, right aligned
$styleDonor15 = "F15"   # numeric style, #,##0
$styleDonor16 = "L15"   # numeric style, #,##0.0 w/ custom negative

function Set-NumberCell($ref, $value) {
    $ws.Range($ref).Value = $value
}

function Set-NumberCellWithStyle($ref, $value, $donor) {
    $ws.Range($ref).Value = $value
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

function Set-StringCell($ref, $value, $donor) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $value
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

# Row 15
Set-StringCell "C15" "0" $styleDonor14
Set-NumberCellWithStyle "D15" 1 $styleDonor15
Set-NumberCellWithStyle "E15" -100 $styleDonor16
Set-NumberCell "F15" 2
Set-NumberCell "G15" 2
Set-NumberCell "H15" 0
Set-NumberCell "I15" 5
Set-NumberCell "J15" 5
Set-NumberCell "K15" 0
Set-NumberCell "L15" 400
Set-NumberCell "M15" 66.666666666666
Set-NumberCell "N15" -44.444444444444

# Row 16
Set-NumberCell "C16" 2
Set-StringCell "D16" "0" $styleDonor14
Set-StringCell "E16" "***.*" $styleDonor14
Set-NumberCell "F16" 12
Set-NumberCell "G16" 9
Set-NumberCell "H16" 33.333333333333
Set-NumberCell "I16" 43
Set-NumberCell "J16" 30
Set-NumberCell "K16" 43.333333333333
Set-NumberCell "L16" 10.25641025641
Set-NumberCell "M16" -21.818181818181
Set-NumberCell "N16" -78.712871287128

# Row 17
Set-NumberCell "C17" 4
Set-NumberCell "D17" 4
Set-NumberCell "E17" 0
Set-NumberCell "F17" 9
Set-NumberCell "G17" 16
Set-NumberCell "H17" -43.75
Set-NumberCell "I17" 36
Set-NumberCell "J17" 38
Set-NumberCell "K17" -5.263157894736
Set-NumberCell "L17" 5.882352941176
Set-NumberCell "M17" 20
Set-NumberCell "N17" -63.636363636363

# Row 18
Set-NumberCellWithStyle "C18" 3 $styleDonor15
Set-NumberCell "D18" 4
Set-NumberCell "E18" -25
Set-NumberCell "F18" 4
Set-NumberCell "G18" 19
Set-NumberCell "H18" -78.947368421052
Set-NumberCell "I18" 32
Set-NumberCell "J18" 41
Set-NumberCell "K18" -21.951219512195
Set-NumberCell "L18" -52.238805970149
Set-NumberCell "M18" -5.882352941176
Set-NumberCell "N18" -89.57654723127

# Row 19
Set-NumberCell "C19" 9
Set-NumberCell "D19" 3
Set-NumberCell "E19" 200
Set-NumberCell "F19" 31
Set-NumberCell "G19" 36
Set-NumberCell "H19" -13.888888888888
Set-NumberCell "I19" 110
Set-NumberCell "J19" 138
Set-NumberCell "K19" -20.289855072463
Set-NumberCell "L19" -9.83606557377
Set-NumberCell "M19" -14.728682170542
Set-NumberCell "N19" -59.107806691449

# Row 20
Set-StringCell "C20" "0" $styleDonor14
Set-NumberCell "D20" 1
Set-NumberCell "E20" -100
Set-NumberCell "F20" 3
Set-NumberCell "G20" 8
Set-NumberCell "H20" -62.5
Set-NumberCell "I20" 13
Set-NumberCell "J20" 29
Set-NumberCell "K20" -55.172413793103
Set-NumberCell "L20" -13.333333333333
Set-NumberCell "M20" 44.444444444444
Set-NumberCell "N20" -94.650205761316

# Row 21
Set-NumberCell "C21" 18
Set-NumberCell "D21" 13
Set-NumberCell "E21" 38.461538461538
Set-NumberCell "F21" 61
Set-NumberCell "G21" 90
Set-NumberCell "H21" -32.222222222222
Set-NumberCell "I21" 239
Set-NumberCell "J21" 283
Set-NumberCell "K21" -15.547703180212
Set-NumberCell "L21" -14.028776978417
Set-NumberCell "M21" -8.429118773946
Set-NumberCell "N21" -78.961267605633

# Row 22
Set-NumberCell "C22" 1
Set-NumberCellWithStyle "D22" 1 $styleDonor15
Set-NumberCellWithStyle "E22" 0 $styleDonor16
Set-NumberCell "F22" 2
Set-NumberCell "G22" 2
Set-NumberCell "H22" 0
Set-NumberCell "I22" 10
Set-NumberCell "J22" 5
Set-NumberCell "K22" 100
Set-NumberCell "L22" 11.111111111111
Set-NumberCell "M22" 11.111111111111

# Row 23
Set-NumberCell "C23" 2
Set-NumberCell "D23" 2
Set-NumberCell "E23" 0
Set-NumberCell "F23" 4
Set-NumberCell "G23" 10
Set-NumberCell "H23" -60
Set-NumberCell "I23" 20
Set-NumberCell "J23" 28
Set-NumberCell "K23" -28.571428571428
Set-NumberCell "L23" -23.076923076923
Set-NumberCell "M23" -28.571428571428

# Row 24
Set-NumberCell "C24" 22
Set-NumberCell "D24" 18
Set-NumberCell "E24" 22.222222222222
Set-NumberCell "F24" 105
Set-NumberCell "G24" 103
Set-NumberCell "H24" 1.941747572815
Set-NumberCell "I24" 343
Set-NumberCell "J24" 485
Set-NumberCell "K24" -29.278350515463
Set-NumberCell "L24" -33.655705996131
Set-NumberCell "M24" 35.03937007874

# Row 25
Set-NumberCell "C25" 13
Set-NumberCell "D25" 8
Set-NumberCell "E25" 62.5
Set-NumberCell "F25" 67
Set-NumberCell "G25" 70
Set-NumberCell "H25" -4.285714285714
Set-NumberCell "I25" 226
Set-NumberCell "J25" 329
Set-NumberCell "K25" -31.306990881459
Set-NumberCell "L25" -40.053050397878

# Row 26
Set-NumberCell "C26" 2
Set-NumberCell "D26" 7
Set-NumberCell "E26" -71.428571428571
Set-NumberCell "F26" 21
Set-NumberCell "G26" 21
Set-NumberCell "H26" 0
Set-NumberCell "I26" 73
Set-NumberCell "J26" 70
Set-NumberCell "K26" 4.285714285714
Set-NumberCell "L26" 7.35294117647
Set-NumberCell "M26" -6.410256410256

# Row 27
Set-StringCell "C27" "0" $styleDonor14
Set-NumberCellWithStyle "D27" 1 $styleDonor15
Set-NumberCellWithStyle "E27" -100 $styleDonor16
Set-NumberCell "F27" 2
Set-NumberCell "G27" 3
Set-NumberCell "H27" -33.333333333333
Set-NumberCell "I27" 5
Set-NumberCell "J27" 7
Set-NumberCell "K27" -28.571428571428

# Row 28
Set-NumberCell "D28" 1
Set-NumberCell "E28" 0
Set-NumberCell "F28" 4
Set-NumberCell "G28" 7
Set-NumberCell "H28" -42.857142857142
Set-NumberCell "I28" 12
Set-NumberCell "J28" 16
Set-NumberCell "K28" -25
Set-NumberCell "L28" -7.692307692307

# Row 31
Set-NumberCell "F31" 5
Set-NumberCell "H31" 25
Set-NumberCell "I31" 8
Set-NumberCell "K31" 100
Set-NumberCell "L31" 33.333333333333

Write-Output "All edits applied"